# Dividend Calculation.xlsx - apply June 2017 dividend correction
# (Yearly!L8:N8 updated; all dependent totals/formulas recalc automatically)

$wb = $excel.ActiveWorkbook

$yearly  = $wb.Worksheets.Item("Yearly")
$allTime = $wb.Worksheets.Item("All Time")

# --- Yearly sheet: corrected June 2017 figures (row 8, 2017 block) ---
$yearly.Range("L8").Value = 118.95   # Taxable Account
$yearly.Range("M8").Value = 44.98    # 401K
$yearly.Range("N8").Value = 42.31    # Suzie's Roth IRA
# O8 = SUM(L8:N8) is a formula and recalculates on its own.

# --- Update the on-screen selection / scroll position for each sheet ---
$yearly.Activate()
$yearly.Range("K22").Select()

$allTime.Activate()
$excel.ActiveWindow.ScrollRow = 25
$allTime.Range("A55").Select()

# Leave "Yearly" as the active (tab-selected) sheet, matching the workbook's
# original active tab.
$yearly.Activate()

$wb.Save()
